# Weekly data refresh: insert this week's new price records (two rows,
# "Primera" and "Segunda" quality) at the top of the data block (row 1117),
# pushing all existing records down by two rows. Excel's native row Insert
# keeps every existing row's content/formatting intact (including the date
# column's number format), which is exactly what the diff shows happened.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new blank rows at 1117-1118; rows 1117..1237 shift down to 1119..1239.
$ws.Rows("1117:1118").Insert()

# Populate the newly inserted row 1117 ("Primera").
$ws.Cells.Item(1117, 1).Value  = 8
$ws.Cells.Item(1117, 2).Value  = "Terminal La Palmera de La Serena"
$ws.Cells.Item(1117, 3).Value  = "Coquimbo"
$ws.Cells.Item(1117, 4).Value  = 45194
$ws.Cells.Item(1117, 5).Value  = 4
$ws.Cells.Item(1117, 6).Value  = 100112023
$ws.Cells.Item(1117, 7).Value  = "Brócoli"
$ws.Cells.Item(1117, 8).Value  = "Sin especificar"
$ws.Cells.Item(1117, 9).Value  = "Primera"
$ws.Cells.Item(1117, 10).Value = 2000
$ws.Cells.Item(1117, 11).Value = 700
$ws.Cells.Item(1117, 12).Value = 800
$ws.Cells.Item(1117, 13).Value = 750
$ws.Cells.Item(1117, 14).Value = "$/unidad"
$ws.Cells.Item(1117, 15).Value = "Provincia del Elquí"
$ws.Cells.Item(1117, 16).Value = 750
$ws.Cells.Item(1117, 17).Value = 1
$ws.Cells.Item(1117, 18).Value = "Hortaliza"

# Populate the newly inserted row 1118 ("Segunda").
$ws.Cells.Item(1118, 1).Value  = 8
$ws.Cells.Item(1118, 2).Value  = "Terminal La Palmera de La Serena"
$ws.Cells.Item(1118, 3).Value  = "Coquimbo"
$ws.Cells.Item(1118, 4).Value  = 45194
$ws.Cells.Item(1118, 5).Value  = 4
$ws.Cells.Item(1118, 6).Value  = 100112023
$ws.Cells.Item(1118, 7).Value  = "Brócoli"
$ws.Cells.Item(1118, 8).Value  = "Sin especificar"
$ws.Cells.Item(1118, 9).Value  = "Segunda"
$ws.Cells.Item(1118, 10).Value = 1160
$ws.Cells.Item(1118, 11).Value = 500
$ws.Cells.Item(1118, 12).Value = 600
$ws.Cells.Item(1118, 13).Value = 550
$ws.Cells.Item(1118, 14).Value = "$/unidad"
$ws.Cells.Item(1118, 15).Value = "Provincia del Elquí"
$ws.Cells.Item(1118, 16).Value = 550
$ws.Cells.Item(1118, 17).Value = 1
$ws.Cells.Item(1118, 18).Value = "Hortaliza"
